$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily-scrape row for 2026/02/06 (金) was inserted above the
# existing 2026/12/29 block, shifting every row from 762..803 down to
# 763..804 (dimension grows from D803 to D804).
$ws.Rows.Item(762).Insert()

# Column A holds date-like text ("2026/02/06"), not a real date value;
# force text formatting first so Excel doesn't auto-convert the string
# into a date serial, then restore the default "Normal" style so the
# new cell doesn't end up with a lingering custom number format.
$ws.Cells.Item(762, 1).NumberFormat = "@"
$ws.Cells.Item(762, 1).Value = "2026/02/06"
$ws.Cells.Item(762, 1).Style = "Normal"

$ws.Cells.Item(762, 2).Value = "金"
$ws.Cells.Item(762, 3).Value = 3
$ws.Cells.Item(762, 4).Value = 55
